$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Datos actualizados" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Agosto de 2020 a las 07:03"

# --- India (row 6) ---
$ws.Range("B6").Value = 3106348
$ws.Range("C6").Value = 1163
$ws.Range("D6").Value = 2338035
$ws.Range("E6").Value = 710621

# --- Pakistan (row 19) ---
$ws.Range("B19").Value = 293261
$ws.Range("C19").Value = 496
$ws.Range("D19").Value = 276829
$ws.Range("E19").Value = 10188
$ws.Range("G19").Value = 9
$ws.Range("H19").Value = 6244

# --- Kirguistan (row 56) ---
$ws.Range("B56").Value = 43126
$ws.Range("C56").Value = 103
$ws.Range("D56").Value = 36615
$ws.Range("E56").Value = 5454
$ws.Range("G56").Value = 1
$ws.Range("H56").Value = 1057

# --- Tailandia / Mozambique swap (rows 120 & 121) ---
# Tailandia overtakes Mozambique in total cases, so it now sits in row 120
# (previously Mozambique's row) with its own updated figures, while
# Mozambique (figures unchanged) drops to row 121.
$ws.Range("A120").Value = "Tailandia"
$ws.Range("B120").Value = 3397
$ws.Range("C120").Value = 2
$ws.Range("D120").Value = 3222
$ws.Range("E120").Value = 117
$ws.Range("H120").Value = 58

$ws.Range("A121").Value = "Mozambique"
$ws.Range("B121").Value = 3395
$ws.Range("C121").Value = 0
$ws.Range("D121").Value = 1503
$ws.Range("E121").Value = 1872
$ws.Range("H121").Value = 20

# --- Butan (row 190) ---
$ws.Range("D190").Value = 115
$ws.Range("E190").Value = 40

# --- Santa Lucia / Timor Oriental swap (rows 202 & 203) ---
# Figures are identical for both countries, only their relative order changes.
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "Santa Lucia"
